# Update the income-statement database: drop the oldest reporting period
# (column D, FY 1396/12) and shift every remaining period one column to
# the left, then append the newly published period (FY 1401/12) in the
# now-empty column H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that hold the period-header label, plus every data row of the
# income statement table: these shift all five columns D:H one period to
# the left (new col = old col + 1), dropping the oldest period.
$rowsFullShift = @(8, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27)

foreach ($r in $rowsFullShift) {
    for ($c = 4; $c -le 7; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $ws.Cells.Item($r, $c + 1).Value2
    }
}

# Row 9 (publish date) only shifts D:F left; the two most recent release
# dates (G, H) are freshly supplied by the new database snapshot rather
# than being carried over from the old H column.
for ($c = 4; $c -le 6; $c++) {
    $ws.Cells.Item(9, $c).Value2 = $ws.Cells.Item(9, $c + 1).Value2
}

# Newly published period: header label
$ws.Cells.Item(8, 8).Value2 = "12 ماهه منتهی به 1401/12"

# Newly published period: financial figures (row -> new column-H value)
$ws.Cells.Item(11, 8).Value2 = 79282474
$ws.Cells.Item(12, 8).Value2 = -25473270
$ws.Cells.Item(13, 8).Value2 = 53809204
$ws.Cells.Item(14, 8).Value2 = -3765171
$ws.Cells.Item(15, 8).Value2 = 0
$ws.Cells.Item(16, 8).Value2 = 3466328
$ws.Cells.Item(17, 8).Value2 = 53510361
$ws.Cells.Item(18, 8).Value2 = -166619
$ws.Cells.Item(19, 8).Value2 = -6761765
$ws.Cells.Item(20, 8).Value2 = 46581977
$ws.Cells.Item(21, 8).Value2 = 0
$ws.Cells.Item(22, 8).Value2 = 46581977
$ws.Cells.Item(23, 8).Value2 = 0
$ws.Cells.Item(24, 8).Value2 = 46581977
$ws.Cells.Item(25, 8).Value2 = 13199
$ws.Cells.Item(26, 8).Value2 = 3529200
$ws.Cells.Item(27, 8).Value2 = 13199

# Publish dates for the two most recent releases.
$ws.Cells.Item(9, 7).Value2 = "1402-02-27 (9)"
$ws.Cells.Item(9, 8).Value2 = "1402-02-27 (2)"
